$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer ID (B1) -- keep leading zeros by formatting the cell as Text first
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "00069790"

# SF Name (F1)
$ws.Range("F1").Value = "Sussie Due"

# ZT01 No. (B3)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "300007240"

# Box Name (F3)
$ws.Range("F3").Value = "Sussie Due"

# VL01N No. (B5) -- keep leading zero
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0880012712"

# Date Recieved (F5)
$ws.Range("F5").Value = "04/23"

# Notes (E10)
$ws.Range("E10").Value = "upon determination, there is a problem as soon as the box is open. please provide warranty service."
